$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply AutoFilter over the original table extent (A1:L62) before the new
# row is appended, so the autoFilter range keeps excluding the new row.
[void]$ws.Range("A1:L62").AutoFilter()

# Update the D5 value (was 0, becomes 776000000).
$ws.Range("D5").Value = 776000000

# Widen column F (6) - closest achievable width to 13.140625 given this
# engine's column-width quantization.
$ws.Columns.Item(6).ColumnWidth = 12.35

# Append the new project row (row 63).
$ws.Range("A63").Value = 9395028
$ws.Range("B63").Value = "Projetos Regionais Rodoviários - Mobilidade regional na Bacia do Paraopeba"
$ws.Range("C63").Value = "I.3"
$ws.Range("D63").Value = 262717753.98
$ws.Range("A63:E63").RowHeight = 30.75

# Shrink the hidden _FilterDatabase defined name down to the header row
# only (matches the post-edit workbook.xml).
$n = $ws.Names("_xlnm._FilterDatabase")
$n.RefersTo = "=projetos!`$A`$1:`$L`$1"

# Move the active selection to E2.
[void]$ws.Range("E2").Select()
